$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.632.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.588.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.810.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.576.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.608.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.680"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +22.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.313.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("E37").Value = "  -5.03%  "
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("E39").Value = "  -3.71%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.792"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.723.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("E48").Value = "  -9.03%  "
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.93%  "
